# Updates cryptos list prices / 1h volume percentages (and re-orders two
# coin rows) to match the latest scrape, per the commit:
#   "Updated cryptos list on Mon Jul  3 09:31:11 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.666.09'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.962.37'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''248.71'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = '''1.0000'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '''0.4825'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '''0.2952'
$ws.Range('E8').Value = '  +2.39%  '
$ws.Range('D9').Value = '''0.06793'
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('D10').Value = '''110.85'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').Value = '''19.38'
$ws.Range('D12').Value = '1.968.41'
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('D13').Value = '''0.07734'
$ws.Range('E13').Value = '  +2.37%  '
$ws.Range('E14').Value = '  +4.69%  '
$ws.Range('D15').Value = '''0.6908'
$ws.Range('E15').Value = '  +3.44%  '
$ws.Range('D16').Value = '''293.82'
$ws.Range('E16').Value = '  -0.64%  '
$ws.Range('D17').Value = '30.680.34'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '''13.29'
$ws.Range('E18').Value = '  +2.58%  '
$ws.Range('D19').Value = '''5.678'
$ws.Range('E19').Value = '  +3.33%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '''0.000007693'
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.218.67'
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = '''0.9999'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '''6.606'
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').Value = '''9.927'
$ws.Range('E25').Value = '  +4.77%  '
$ws.Range('D26').Value = '''170.96'
$ws.Range('E26').Value = '  +3.98%  '
$ws.Range('D27').Value = '''20.14'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('E28').Value = '  +4.49%  '
$ws.Range('D29').Value = '''0.1072'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '''1.440'
$ws.Range('E30').Value = '  +2.73%  '
$ws.Range('D31').Value = '''4.696'
$ws.Range('E31').Value = '  +16.60%  '
$ws.Range('D32').Value = '''4.465'
$ws.Range('E32').Value = '  +7.06%  '
$ws.Range('D33').Value = '''0.05130'
$ws.Range('E33').Value = '  +3.18%  '
$ws.Range('D34').Value = '''0.7809'
$ws.Range('E34').Value = '  +7.02%  '
$ws.Range('D35').Value = '''1.181'
$ws.Range('E35').Value = '  +4.00%  '
$ws.Range('D36').Value = '''0.02064'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').Value = '''2.734'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '''2.710'
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('D39').Value = '''2.066'
$ws.Range('E39').Value = '  +2.58%  '
$ws.Range('D40').Value = '''111.12'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('D41').Value = '''6.132'
$ws.Range('E41').Value = '  +4.43%  '
$ws.Range('D42').Value = '''0.4471'
$ws.Range('E42').Value = '  +1.20%  '
$ws.Range('D43').Value = '''0.8749'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D44').Value = '''70.26'
$ws.Range('E44').Value = '  +2.65%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = '''7.397'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').Value = '''0.1277'
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('D48').Value = '''9.352'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D49').Value = '''35.89'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').Value = '''47.81'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '''0.4091'
$ws.Range('E51').Value = '  +1.93%  '

# The quote-prefix trick above marks the cell with a "number stored as
# text" style (quotePrefix). Reset the style back to Normal so the cell
# keeps its original plain formatting, matching the source workbook.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
